# Regenerate the "K" (strikeouts) column (column G) of the save_data sheet
# so that it reflects actual strikeout counts (K) instead of the old
# "Strike#" pitch-count-based values. Values below were recalculated
# (regen std/mean, calc and write s_vals) from the source stats for this
# game log and are written directly into column G, rows 2-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 3
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 3
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    16 = 2
    17 = 3
    18 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
